$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '64.163.86'
$ws.Cells.Item(2, 5).Value = '  -3.37%  '
$ws.Cells.Item(3, 4).Value = '3.174.72'
$ws.Cells.Item(3, 5).Value = '  -8.30%  '
$ws.Cells.Item(4, 5).Value = '  -0.01%  '
$ws.Cells.Item(5, 4).Value = '564.89'
$ws.Cells.Item(5, 5).Value = '  -3.22%  '
$ws.Cells.Item(6, 4).Value = '168.97'
$ws.Cells.Item(6, 5).Value = '  -5.00%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.610'
$ws.Cells.Item(7, 4).ClearFormats()
$ws.Cells.Item(7, 5).Value = '  -3.01%  '
$ws.Cells.Item(8, 5).Value = '  +0.02%  '
$ws.Cells.Item(9, 4).Value = '3.171.98'
$ws.Cells.Item(9, 5).Value = '  -8.33%  '
$ws.Cells.Item(10, 5).Value = '  -6.87%  '
$ws.Cells.Item(11, 4).Value = '6.59'
$ws.Cells.Item(11, 5).Value = '  -5.33%  '
$ws.Cells.Item(12, 5).Value = '  -5.43%  '
$ws.Cells.Item(13, 4).Value = '3.721.64'
$ws.Cells.Item(13, 5).Value = '  -8.38%  '
$ws.Cells.Item(14, 5).Value = '  +1.39%  '
$ws.Cells.Item(15, 5).Value = '  -9.44%  '
$ws.Cells.Item(16, 4).Value = '64.139.95'
$ws.Cells.Item(16, 5).Value = '  -3.24%  '
$ws.Cells.Item(17, 5).Value = '  -5.44%  '
$ws.Cells.Item(18, 4).Value = '3.172.85'
$ws.Cells.Item(18, 5).Value = '  -8.03%  '
$ws.Cells.Item(19, 4).Value = '5.72'
$ws.Cells.Item(19, 5).Value = '  -4.42%  '
$ws.Cells.Item(20, 5).Value = '  -6.43%  '
$ws.Cells.Item(21, 4).Value = '352.51'
$ws.Cells.Item(21, 5).Value = '  -5.09%  '
$ws.Cells.Item(22, 5).Value = '  -6.38%  '
$ws.Cells.Item(23, 5).Value = '  +0.44%  '
$ws.Cells.Item(24, 4).Value = '68.56'
$ws.Cells.Item(25, 4).Value = '0.505'
$ws.Cells.Item(25, 5).Value = '  -5.95%  '
$ws.Cells.Item(26, 5).Value = '  -6.03%  '
$ws.Cells.Item(27, 5).Value = '  -4.34%  '
$ws.Cells.Item(28, 5).Value = '  -1.00%  '
$ws.Cells.Item(29, 5).Value = '  -0.15%  '
$ws.Cells.Item(30, 5).Value = '  -0.17%  '
$ws.Cells.Item(31, 4).Value = '5.56'
$ws.Cells.Item(31, 5).Value = '  -6.97%  '
$ws.Cells.Item(32, 5).Value = '  -4.97%  '
$ws.Cells.Item(33, 4).Value = '21.95'
$ws.Cells.Item(34, 5).Value = '  -6.22%  '
$ws.Cells.Item(35, 5).Value = '  -5.39%  '
$ws.Cells.Item(36, 5).Value = '  -8.29%  '
$ws.Cells.Item(37, 4).Value = '153.86'
$ws.Cells.Item(37, 5).Value = '  -4.60%  '
$ws.Cells.Item(38, 4).Value = '0.817'
$ws.Cells.Item(38, 5).Value = '  -7.89%  '
$ws.Cells.Item(39, 4).Value = '26.02'
$ws.Cells.Item(39, 5).Value = '  -6.82%  '
$ws.Cells.Item(40, 5).Value = '  -6.49%  '
$ws.Cells.Item(41, 5).Value = '  -4.28%  '
$ws.Cells.Item(42, 4).Value = '2.613.64'
$ws.Cells.Item(42, 5).Value = '  -7.03%  '
$ws.Cells.Item(43, 5).Value = '  -7.74%  '
$ws.Cells.Item(45, 4).Value = '39.28'
$ws.Cells.Item(45, 5).Value = '  -1.89%  '
$ws.Cells.Item(46, 4).Value = '0.0645'
$ws.Cells.Item(46, 5).Value = '  -6.91%  '
$ws.Cells.Item(47, 4).Value = '23.65'
$ws.Cells.Item(47, 5).Value = '  -6.60%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '317.60'
$ws.Cells.Item(48, 4).ClearFormats()
$ws.Cells.Item(48, 5).Value = '  -7.09%  '
$ws.Cells.Item(49, 5).Value = '  -7.14%  '
$ws.Cells.Item(50, 5).Value = '  -3.56%  '
$ws.Cells.Item(51, 5).Value = '  -0.09%  '
